$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update probability values in column C ---
$ws.Range("C4").Value = 0.15
$ws.Range("C5").Value = 0.2
$ws.Range("C8").Value = 0.4

# --- Update text values ---
# Row 5: measure for "Robo de equipos..." changes to the access control text
$ws.Range("D5").Value = "El acceso al área de servidores o a la habitación del rack está restringido al personal del Área de Sistemas y Calidad, los cuales, para desbloquear la correspondiente puerta, se identifican mediante su huella dactilar. Vigilancia las 24 horas del día y alarmas de monitoreo."

# Row 7: Risk becomes "Incendio", measure becomes extinguishers text (set first so the new
# shared string for this text is appended before the others, matching the target order)
$ws.Range("A7").Value = "Incendio"
$ws.Range("D7").Value = "Se cuentan con extintores y sistemas de irrigación (con detectores de humo)."

# Row 6: Risk becomes the detailed flooding description, measure becomes the detailed flood-protection text
$ws.Range("D6").Value = "Los servidores estarán ubicados a 1 m. de altura, a salvo de posibles inundaciones. Además la sala cuenta con alcantarillas para desagotar rápidamente cualquier fuga de agua."
$ws.Range("A6").Value = "Inundación de la sala de servidores causada por daños en las cañerías del baño cercano (o baños de pisos superiores)"
$ws.Rows.Item(6).RowHeight = 60

# Row 8: Risk becomes "Corte de energía eléctrica", measure becomes the UPS text, row height grows
$ws.Range("A8").Value = "Corte de energía eléctrica"
$ws.Range("D8").Value = "Para cortes de energía de corto tiempo, hay asignadas UPSs para mantener la alimentación de los servidores. Además, se cuenta con generadores eléctricos para cortes prolongados de más de 1 hora."
$ws.Rows.Item(8).RowHeight = 78.75

# Apply the new style (border + centered/wrapped alignment) to A6, matching the workbook's added cellXf
$ws.Range("A6").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A6").VerticalAlignment = -4108    # xlCenter
$ws.Range("A6").WrapText = $true
$ws.Range("A6").Borders.LineStyle = 1

# --- Update the sheet view / selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("D8").Select()

$wb.Save()
